# Updates the Coin / Link / Price / Volume(1h) table on Sheet1 (rows 2-51)
# to match the refreshed cryptos.xlsx snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.243.95"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3
$ws.Range("D3").Value = "1.862.57"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'235.85"

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").Value = "'0.4672"
$ws.Range("E7").Value = "  -0.55%  "

# Row 8
$ws.Range("D8").Value = "'0.2833"
$ws.Range("E8").Value = "  +0.60%  "

# Row 9
$ws.Range("D9").Value = "'0.06515"
$ws.Range("E9").Value = "  -0.57%  "

# Row 10
$ws.Range("D10").Value = "'21.45"
$ws.Range("E10").Value = "  +6.65%  "

# Row 11
$ws.Range("D11").Value = "'0.07909"
$ws.Range("E11").Value = "  +1.25%  "

# Row 12
$ws.Range("D12").Value = "'97.08"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
$ws.Range("D13").Value = "1.862.57"
$ws.Range("E13").Value = "  +0.15%  "

# Row 14
$ws.Range("D14").Value = "'5.154"
$ws.Range("E14").Value = "  +0.98%  "

# Row 15
$ws.Range("D15").Value = "'0.6782"
$ws.Range("E15").Value = "  +1.94%  "

# Row 16
$ws.Range("D16").Value = "'278.23"
$ws.Range("E16").Value = "  -1.92%  "

# Row 17
$ws.Range("D17").Value = "30.241.97"
$ws.Range("E17").Value = "  -0.07%  "

# Row 18
$ws.Range("D18").Value = "'13.71"
$ws.Range("E18").Value = "  +8.95%  "

# Row 19
$ws.Range("E19").Value = "  +0.14%  "

# Row 20
$ws.Range("D20").Value = "'5.384"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.000007302"
$ws.Range("E21").Value = "  +0.72%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.106.52"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").Value = "'6.148"
$ws.Range("E24").Value = "  +0.14%  "

# Row 25
$ws.Range("D25").Value = "'167.26"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26
$ws.Range("D26").Value = "'9.154"
$ws.Range("E26").Value = "  -1.70%  "

# Row 27
$ws.Range("D27").Value = "'19.02"
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("D28").Value = "'1.923"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29
$ws.Range("D29").Value = "'1.385"
$ws.Range("E29").Value = "  +3.43%  "

# Row 30
$ws.Range("D30").Value = "'0.09725"
$ws.Range("E30").Value = "  +1.21%  "

# Row 31
$ws.Range("D31").Value = "'4.364"
$ws.Range("E31").Value = "  -1.06%  "

# Row 32
$ws.Range("D32").Value = "'1.479"
$ws.Range("E32").Value = "  +0.42%  "

# Row 33
$ws.Range("D33").Value = "'4.023"
$ws.Range("E33").Value = "  -1.91%  "

# Row 34
$ws.Range("E34").Value = "  +0.94%  "

# Row 35
$ws.Range("D35").Value = "'1.126"
$ws.Range("E35").Value = "  +2.23%  "

# Row 36
$ws.Range("D36").Value = "'0.7040"
$ws.Range("E36").Value = "  +0.54%  "

# Row 37
$ws.Range("D37").Value = "'2.707"
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
$ws.Range("D38").Value = "'0.01858"
$ws.Range("E38").Value = "  +0.19%  "

# Row 39
$ws.Range("D39").Value = "'2.585"
$ws.Range("E39").Value = "  +2.98%  "

# Row 40
$ws.Range("D40").Value = "'6.322"
$ws.Range("E40").Value = "  -1.80%  "

# Row 41
$ws.Range("D41").Value = "'75.01"
$ws.Range("E41").Value = "  +4.08%  "

# Row 42
$ws.Range("D42").Value = "'1.956"
$ws.Range("E42").Value = "  +0.89%  "

# Row 43
$ws.Range("D43").Value = "'0.8490"
$ws.Range("E43").Value = "  -0.81%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4166"
$ws.Range("E44").Value = "  +0.12%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9998"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").Value = "'103.44"
$ws.Range("E46").Value = "  -0.57%  "

# Row 47
$ws.Range("D47").Value = "'978.98"
$ws.Range("E47").Value = "  -2.59%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.316"
$ws.Range("E48").Value = "  +3.56%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.154"
$ws.Range("E49").Value = "  -0.70%  "

# Row 50
$ws.Range("D50").Value = "'34.01"
$ws.Range("E50").Value = "  +0.76%  "

# Row 51
$ws.Range("D51").Value = "'0.05648"
$ws.Range("E51").Value = "  +0.06%  "
